# Remove the trailing "Hands-on Session (CNN)" slide (slide 16, the
# last slide, laid out with the unused "Big number" custom layout) along
# with its notes page. Deleting the slide also drops its notes slide
# (notesSlide16.xml) since nothing else references it.
$p = $ppt.ActivePresentation
$p.Slides.Item(16).Delete()

# The "Big number" layout (slideLayout7.xml / CustomLayouts index 7) was
# only used by that slide; remove the now-unused layout from the master
# too, matching the target deck.
$p.SlideMaster.CustomLayouts.Item(7).Delete()
